$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Equity Portfolio
$ws.Range("B2").Value = 127.1376995041531
$ws.Range("C2").Value = 14.76691603482945
$ws.Range("D2").Value = 24.05913949500928
$ws.Range("E2").Value = 0.5306472427028307
$ws.Range("F2").Value = -34.39790091985317

# Row 3 - Covered Call Strategy
$ws.Range("B3").Value = 182.3086931523248
$ws.Range("C3").Value = 19.03405453747926
$ws.Range("D3").Value = 19.38685065985675
$ws.Range("E3").Value = 0.8786395911508574
$ws.Range("F3").Value = -15.18092499036337

# Row 4 - Combined Portfolio
$ws.Range("B4").Value = 151.5557324430735
$ws.Range("C4").Value = 16.75129530213459
$ws.Range("D4").Value = 19.87258352320957
$ws.Range("E4").Value = 0.7422937880676802
$ws.Range("F4").Value = -24.98145958204303

# Row 5 - SPY Buy & Hold
$ws.Range("B5").Value = 125.1811750865068
$ws.Range("C5").Value = 14.60034692120322
$ws.Range("D5").Value = 19.70244516213175
$ws.Range("E5").Value = 0.6395321401742148
$ws.Range("F5").Value = -33.71725745991026
